$d = $word.ActiveDocument
$pairs = @(
    ,@('2025-06-26', '2025-06-30')
    ,@('Date: 2025-06-26', 'Date: 2025-06-30')
    ,@('risk, hidden costs', 'Risk, hidden costs')
    ,@('principles apply to ANY software', 'Principles apply to ANY software')
    ,@('risk management', 'Risk management')
    ,@('due diligence (maintenance, updates)', 'Due diligence (maintenance, updates)')
    ,@('lifecycle management/EOL → migration', 'Lifecycle management/EOL → migration')
    ,@('total cost of ownership considerations', 'Total cost of ownership considerations')
    ,@('benefits (freedom, cost, reducing vendor lock in, portability)', 'Benefits (freedom, cost, reducing vendor lock in, portability)')
    ,@('infrastructure considerations', 'Infrastructure considerations')
    ,@('national policies', 'National policies')
    ,@('events/hackathons (eg. OGC/OSGeo/ASF Joint Sprints)', 'Events/hackathons (eg. OGC/OSGeo/ASF Joint Sprints)')
    ,@('by product: connection/collab', 'By product: connection/collab')
    ,@('regulations / risk / constraints / considerations', 'Regulations / risk / constraints / considerations')
    ,@('aligning with WMO standards', 'Aligning with WMO standards')
    ,@('achieving compliance', 'Achieving compliance')
    ,@('coordination/support functions', 'Coordination/support functions')
    ,@('software selection for WMO application development', 'Software selection for WMO projects and application development')
    ,@('managing FOSS activities', 'Managing FOSS activities')
    ,@('ensuring sustainability of FOSS usage', 'Ensuring sustainability of FOSS usage')
    ,@('managing risk', 'Managing risk')
    ,@('functions', 'Functions')
    ,@('people', 'People')
    ,@('compatability / compliance matrix', 'Compatability / compliance matrix')
    ,@('implementation of WMO Tech Regs / compliance ?', 'Implementation of WMO Tech Regs / compliance ?')
    ,@('ensure FOSS implementations are part of Technical Regulation development/assessment (feasibility)', 'Ensure FOSS implementations are part of Technical Regulation development/assessment (feasibility)')
    ,@('example: wis2box, developed at the same time as WIS2 standards', 'Example: wis2box, developed at the same time as WIS2 standards')
    ,@('example: OGC standards (3 implementations)', 'Example: OGC standards (3 implementations)')
    ,@('software identification and selection', 'Software identification and selection')
    ,@('project checklist/assessment', 'Project checklist/assessment')
    ,@('"approved projects" and/or Reference Implementations', '"Approved projects" and/or Reference Implementations')
    ,@('make Tech Regs more concrete', 'Make Tech Regs more concrete')
    ,@('should FOSS be cited in WMO Tech Regs (suggest no)', 'Should FOSS be cited in WMO Tech Regs (suggest no)')
    ,@('criteria needed', 'Criteria needed')
    ,@('compliance (data exchange)', 'Compliance (data exchange)')
    ,@('software evaluation (FOSS!) checklist → confidence', 'Software evaluation (FOSS!) checklist → confidence')
    ,@('readiness', 'Readiness')
    ,@('bus factor', 'Bus/retirement factor')
    ,@('rolling review', 'Rolling review')
    ,@('harmonization: regular review of ecosystem to ensure alignment and optimal use of resources', 'Harmonization: regular review of ecosystem to ensure alignment and optimal use of resources')
    ,@('case study: wis2box et. al.', 'Case study: wis2box et. al.')
    ,@('agile development during Tech Reg development', 'Agile development during Tech Reg development')
)

$applied = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -eq 0) { continue }
    $clean = $t.TrimEnd([char]13, [char]7)
    foreach ($pair in $pairs) {
        if ($clean.Equals($pair[0], [System.StringComparison]::Ordinal)) {
            $p.Range.Text = $pair[1]
            $applied = $applied + 1
            break
        }
    }
}
Write-Host "Applied: $applied of $($pairs.Count)"
